# Auto-generated edit script applying the Ragnarok_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 1425.4
$ws.Range("I53").Value = 1406.3334
$ws.Range("K53").Value = 1406.3334
$ws.Range("M53").Value = -769.3334

$ws.Range("H114").Value = 97499.5
$ws.Range("J114").Value = 97499.5
$ws.Range("L114").Value = 97499.5
$ws.Range("N114").Value = -106177.5

$ws.Range("H117").Value = 66499
$ws.Range("J117").Value = 66499
$ws.Range("L117").Value = 66499
$ws.Range("N117").Value = -75677

$ws.Range("H120").Value = 99999
$ws.Range("J120").Value = 99999
$ws.Range("L120").Value = 99999
$ws.Range("N120").Value = -109675

$ws.Range("H123").Value = 84999.94
$ws.Range("J123").Value = 84999.94
$ws.Range("L123").Value = 84999.94
$ws.Range("N123").Value = -94799.94

$ws.Range("H124").Value = 99891.5
$ws.Range("J124").Value = 99891.5
$ws.Range("L124").Value = 99891.5
$ws.Range("N124").Value = -109711.5

$ws.Range("H126").Value = 99998.5
$ws.Range("J126").Value = 99998.5
$ws.Range("L126").Value = 99998.5
$ws.Range("N126").Value = -109878.5

$ws.Range("H137").Value = 15153908
$ws.Range("I137").Value = 31251904
$ws.Range("J137").Value = 2852.1177
$ws.Range("K137").Value = 93755712
$ws.Range("L137").Value = 8556.3531
$ws.Range("M137").Value = -93753162
$ws.Range("N137").Value = -13656.3531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H18").Value = 2500
$ws.Range("J18").Value = 2500
$ws.Range("L18").Value = 2500
$ws.Range("N18").Value = -3144

$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()

$ws.Range("H61").Value = 13152574
$ws.Range("I61").Value = 14897730
$ws.Range("J61").Value = 1435099.9
$ws.Range("K61").Value = 14897730
$ws.Range("L61").Value = 1435099.9
$ws.Range("M61").Value = -14897518
$ws.Range("N61").Value = -1435523.9

$ws.Range("H110").Value = 4968.1333
$ws.Range("I110").Value = 4873.88
$ws.Range("K110").Value = 4873.88
$ws.Range("M110").Value = -2828.88

$ws.Range("H136").Value = 13152574
$ws.Range("I136").Value = 14897730
$ws.Range("J136").Value = 1435099.9
$ws.Range("K136").Value = 44693190
$ws.Range("L136").Value = 4305299.699999999
$ws.Range("M136").Value = -44690640
$ws.Range("N136").Value = -4310399.699999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 42886.47
$ws.Range("I86").Value = 69246.7
$ws.Range("J86").Value = 5229
$ws.Range("K86").Value = 69246.7
$ws.Range("L86").Value = 5229
$ws.Range("M86").Value = -68123.7
$ws.Range("N86").Value = -7475

$ws.Range("H89").Value = 42886.47
$ws.Range("I89").Value = 69246.7
$ws.Range("J89").Value = 5229
$ws.Range("K89").Value = 346233.5
$ws.Range("L89").Value = 26145
$ws.Range("M89").Value = -340617.5
$ws.Range("N89").Value = -37377

$ws.Range("H105").Value = 1145904
$ws.Range("I105").Value = 1758968.1
$ws.Range("J105").Value = 7356.4287
$ws.Range("K105").Value = 1758968.1
$ws.Range("L105").Value = 7356.4287
$ws.Range("M105").Value = -1757221.1
$ws.Range("N105").Value = -10850.4287

$ws.Range("H134").Value = 3127035.2
$ws.Range("I134").Value = 2063.04
$ws.Range("J134").Value = 14287650
$ws.Range("K134").Value = 6189.12
$ws.Range("L134").Value = 42862950
$ws.Range("M134").Value = -3654.12
$ws.Range("N134").Value = -42868020

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3119.4
$ws.Range("I58").Value = 3198
$ws.Range("J58").Value = 3001.5
$ws.Range("K58").Value = 3198
$ws.Range("L58").Value = 3001.5
$ws.Range("M58").Value = -2995
$ws.Range("N58").Value = -3407.5

$ws.Range("H86").Value = 13455.818
$ws.Range("I86").Value = 16856
$ws.Range("K86").Value = 16856
$ws.Range("M86").Value = -15733

$ws.Range("H89").Value = 13455.818
$ws.Range("I89").Value = 16856
$ws.Range("K89").Value = 84280
$ws.Range("M89").Value = -78664

$ws.Range("H132").Value = 4408.357
$ws.Range("I132").Value = 2577.5
$ws.Range("J132").Value = 15393.5
$ws.Range("K132").Value = 7732.5
$ws.Range("L132").Value = 46180.5
$ws.Range("M132").Value = -5202.5
$ws.Range("N132").Value = -51240.5

$ws.Range("H134").Value = 2186.5652
$ws.Range("I134").Value = 2437.5625
$ws.Range("K134").Value = 7312.6875
$ws.Range("M134").Value = -4777.6875

$ws.Range("H136").Value = 3119.4
$ws.Range("I136").Value = 3198
$ws.Range("J136").Value = 3001.5
$ws.Range("K136").Value = 9594
$ws.Range("L136").Value = 9004.5
$ws.Range("M136").Value = -7044
$ws.Range("N136").Value = -14104.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 36
$ws.Range("I38").Value = 36.75
$ws.Range("K38").Value = 110.25
$ws.Range("M38").Value = 236.75

$ws.Range("H63").Value = 20490.117
$ws.Range("I63").Value = 5011.5
$ws.Range("J63").Value = 22553.934
$ws.Range("K63").Value = 15034.5
$ws.Range("L63").Value = 67661.802
$ws.Range("M63").Value = -14285.5
$ws.Range("N63").Value = -69159.802

$ws.Range("H66").Value = 20490.117
$ws.Range("I66").Value = 5011.5
$ws.Range("J66").Value = 22553.934
$ws.Range("K66").Value = 45103.5
$ws.Range("L66").Value = 202985.406
$ws.Range("M66").Value = -41359.5
$ws.Range("N66").Value = -210473.406

$ws.Range("H113").Value = 2209.182
$ws.Range("J113").Value = 2899.6
$ws.Range("L113").Value = 8698.799999999999
$ws.Range("N113").Value = -13038.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 409199.4
$ws.Range("I3").Value = 1999999
$ws.Range("J3").Value = 11499.5
$ws.Range("K3").Value = 1999999
$ws.Range("L3").Value = 11499.5
$ws.Range("M3").Value = -1999883
$ws.Range("N3").Value = -11731.5

$ws.Range("H10").Value = 7999.5
$ws.Range("I10").Value = 799
$ws.Range("J10").Value = 10399.667
$ws.Range("K10").Value = 799
$ws.Range("L10").Value = 10399.667
$ws.Range("M10").Value = -630
$ws.Range("N10").Value = -10737.667

$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()

$ws.Range("H14").Value = 3966.6667
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 950
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 950
$ws.Range("M14").Value = -9832
$ws.Range("N14").Value = -1286

$ws.Range("H19").Value = 5005000
$ws.Range("I19").Value = 5005000
$ws.Range("K19").Value = 5005000
$ws.Range("M19").Value = -5004712

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1846.8334
$ws.Range("I82").Value = 1215.6
$ws.Range("K82").Value = 1215.6
$ws.Range("M82").Value = -854.5999999999999

$ws.Range("H85").Value = 1846.8334
$ws.Range("I85").Value = 1215.6
$ws.Range("K85").Value = 1215.6
$ws.Range("M85").Value = 32.40000000000009

$ws.Range("H122").Value = 3848.8298
$ws.Range("I122").Value = 3452.5789
$ws.Range("J122").Value = 5521.8887
$ws.Range("K122").Value = 10357.7367
$ws.Range("L122").Value = 16565.6661
$ws.Range("M122").Value = -7907.736699999999
$ws.Range("N122").Value = -21465.6661

$ws.Range("H132").Value = 3223.7334
$ws.Range("I132").Value = 1949.7646
$ws.Range("K132").Value = 5849.293799999999
$ws.Range("M132").Value = -3319.293799999999

$ws.Range("H136").Value = 4356.92
$ws.Range("I136").Value = 4569.8335
$ws.Range("K136").Value = 13709.5005
$ws.Range("M136").Value = -11159.5005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2149.1765
$ws.Range("I122").Value = 2328.3044
$ws.Range("K122").Value = 6984.9132
$ws.Range("M122").Value = -4534.9132

$ws.Range("H126").Value = 3050.1667
$ws.Range("I126").Value = 3730.5
$ws.Range("J126").Value = 2029.6666
$ws.Range("K126").Value = 11191.5
$ws.Range("L126").Value = 6088.9998
$ws.Range("M126").Value = -8721.5
$ws.Range("N126").Value = -11028.9998

$ws.Range("H132").Value = 1252172.2
$ws.Range("I132").Value = 2482.5715
$ws.Range("K132").Value = 7447.7145
$ws.Range("M132").Value = -4917.7145

$ws.Range("H136").Value = 477614.34
$ws.Range("I136").Value = 1901.6364
$ws.Range("J136").Value = 1000898.3
$ws.Range("K136").Value = 5704.9092
$ws.Range("L136").Value = 3002694.9
$ws.Range("M136").Value = -3154.9092
$ws.Range("N136").Value = -3007794.9
